# Generate Report for Handback
# - Marks rows as handed back (status text updated)
# - Fills in the "Latest Target File" / "Latest Handback File" columns (F/G)
#   with hyperlinks, for both data rows, on the zh-cn and de-de sheets
# - Updates the "Latest Handback DateTime" column (H) with real timestamps

$wb = $excel.ActiveWorkbook

$status = "Handed back: in sync with en-US"

# ---------------------------------------------------------------------
# Overview sheet: refresh the per-language status columns (B, C)
# ---------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("B2").Value = $status
$wsOverview.Range("C2").Value = $status
$wsOverview.Range("B3").Value = $status
$wsOverview.Range("C3").Value = $status

# ---------------------------------------------------------------------
# zh-cn sheet
# ---------------------------------------------------------------------
$wsZh = $wb.Worksheets.Item("zh-cn")

$wsZh.Range("C2").Value = $status
$wsZh.Range("C3").Value = $status

$zhMdUrl = "https://github.com/OpenLocalizationTest/oltest/blob/108e1090d0794f6bc17451c43364088dd1758322/e2e/a.md"
$zhXlfUrl = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/422b50eb09acb845a8102fbca7ada86b564c7a72/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.zh-cn.xlf"
$zhXlfName = "a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.zh-cn.xlf"

$wsZh.Hyperlinks.Add($wsZh.Range("F2"), $zhMdUrl, "", "", "a.md")
$wsZh.Hyperlinks.Add($wsZh.Range("G2"), $zhXlfUrl, "", "", $zhXlfName)
$wsZh.Hyperlinks.Add($wsZh.Range("F3"), $zhMdUrl, "", "", "a.md")
$wsZh.Hyperlinks.Add($wsZh.Range("G3"), $zhXlfUrl, "", "", $zhXlfName)

$wsZh.Range("H2").Value = "2016-03-18 12:23:44"
$wsZh.Range("H3").Value = "2016-03-18 12:23:44"

# ---------------------------------------------------------------------
# de-de sheet
# ---------------------------------------------------------------------
$wsDe = $wb.Worksheets.Item("de-de")

$wsDe.Range("C2").Value = $status
$wsDe.Range("C3").Value = $status

$deMdUrl = "https://github.com/OpenLocalizationTest/oltest/blob/108e1090d0794f6bc17451c43364088dd1758322/e2e/a.md"
$deXlfUrl = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/0b9f80ae6875ebb883bb2f5ba85e1e75dcb73e27/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.de-de.xlf"
$deXlfName = "a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.de-de.xlf"

$wsDe.Hyperlinks.Add($wsDe.Range("F2"), $deMdUrl, "", "", "a.md")
$wsDe.Hyperlinks.Add($wsDe.Range("G2"), $deXlfUrl, "", "", $deXlfName)
$wsDe.Hyperlinks.Add($wsDe.Range("F3"), $deMdUrl, "", "", "a.md")
$wsDe.Hyperlinks.Add($wsDe.Range("G3"), $deXlfUrl, "", "", $deXlfName)

$wsDe.Range("H2").Value = "2016-03-18 12:23:50"
$wsDe.Range("H3").Value = "2016-03-18 12:23:50"
